$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "all": add the 2020-05-05 (serial 43956) daily row just
# above the footnote row (old row 28 -> new row 29).
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")
$wsAll.Rows.Item(28).Insert()

$wsAll.Cells.Item(28, 1).Value = 43956
$wsAll.Cells.Item(28, 2).Value = 269
$wsAll.Cells.Item(28, 3).Value = 263
$wsAll.Cells.Item(28, 4).Value = 121
$wsAll.Cells.Item(28, 5).Value = 112
$wsAll.Cells.Item(28, 6).Value = 9
$wsAll.Cells.Item(28, 7).Value = 7
$wsAll.Cells.Item(28, 8).Value = 135

$wsAll.Activate()
$wsAll.Range("B31").Select()

# ---------------------------------------------------------------
# Sheet "kobe": update the running total for 5/4 and add the
# 2020-05-05 daily row just above the footnote row
# (old row 83 -> new row 84).
# ---------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Cells.Item(82, 4).Value = 2

$wsKobe.Rows.Item(83).Insert()

$wsKobe.Cells.Item(83, 1).Value = 43956
$wsKobe.Cells.Item(83, 3).Value = 2242
$wsKobe.Cells.Item(83, 4).Value = 1
$wsKobe.Cells.Item(83, 5).Value = 269
$wsKobe.Cells.Item(83, 6).Value = 116
$wsKobe.Cells.Item(83, 7).Value = 108
$wsKobe.Cells.Item(83, 8).Value = 8
$wsKobe.Cells.Item(83, 9).Value = 7
$wsKobe.Cells.Item(83, 10).Value = 128

$wsKobe.Activate()
$wsKobe.Range("G84").Select()

# ---------------------------------------------------------------
# Sheet "other": add the 2020-05-05 daily row just above the
# footnote row (old row 58 -> new row 59).
# ---------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows.Item(58).Insert()

$wsOther.Cells.Item(58, 1).Value = 43956
$wsOther.Cells.Item(58, 2).Value = 0
$wsOther.Cells.Item(58, 3).Value = 12
$wsOther.Cells.Item(58, 4).Value = 5
$wsOther.Cells.Item(58, 5).Value = 4
$wsOther.Cells.Item(58, 6).Value = 1
$wsOther.Cells.Item(58, 7).Value = 0
$wsOther.Cells.Item(58, 8).Value = 7

$wsOther.Activate()
$wsOther.Range("G59").Select()

$wsAll.Activate()
